$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the numeric-looking identifier columns (D/E/F) as Text so that
# values like "132" or "52215867" are stored as strings, not numbers -
# matching the source data (t="inlineStr") instead of Excel's default numeric
# auto-detection.
$ws.Range("D2:F4").NumberFormat = "@"

# Row 2
$ws.Range("A2").Value = "281474991205262-1738710055344"
$ws.Range("B2").Value = "Mobile Usage"
$ws.Range("C2").Value = "2025-02-04T17:00:55.344"
$ws.Range("D2").Value = "281474991205262"
$ws.Range("E2").Value = "132"
$ws.Range("F2").Value = "52215867"
$ws.Range("G2").Value = "EMMANUEL SALCEDO"
$ws.Range("H2").Value = 20.63466236
$ws.Range("I2").Value = -103.33006352
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = "https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1738710052844/FbWwjUaCGj-camera-video-segment-driver-1738710055344.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSK5OXEILC%2F20250205%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250205T150543Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjECsaCXVzLXdlc3QtMiJHMEUCIHsgChGSRKwjEjGoDqjcGh0LrGtLRUqAEs2OolKwO3gIAiEAp4o6%2BVCBhl1VcA%2FXVT6kmUoc4Q%2BBKuNSqGQ4TpPl2m4q3QMIRBAEGgw3ODEyMDQ5NDIyNDQiDI4RT2evNNQhplb4uiq6A9490COmLc8PeJwtKgtCKCn4jpIrTNdV5vXkE6bNQZQffREVTswCw8bTKS3pu5bBuIhAD6V8hA20CSqTdb%2BRUm8yG1Zrie012UZIgvDGLl7YkachS6vz4Yx35THvqPLu%2B9Gmb4gHx3IZW7a9piMFic8LVgMcuPVC647r6O5Ip8SIguX3TDbP0fUliyZB8NqvosYJUBjV5nScUPyAQSz0uR8a9%2FjIOBQbXgiYz07uDubW2VC1aYQoRWQks8T0bMeHny748fSDpEBtnAbGly4BXSMRzR%2BNbIy4CHz1CZYruHs%2BWV5SWee6BPBegyYm1fO2aMY1upjqIn89nbNIFtBkgmpCXRn2H4XAUe1zwU%2FSsptShLhmjIhtaU%2BrIw1CGPirdG5%2FVNxN3NUAMoajto%2Fw6RNKdvA8tspWKn1m5InxeYN164ODNASIzkJE5dHFIX%2B7z%2BJzKC0zbUvXVpQKZbYD6vL6nYoEXBXuI6ZnPetxTIOYfmECP1IlJWR59HLdsCboc0zrS7bz%2BAf5AAaTt14Np4KZKFR9G0Phzo4j9EVle1yi7E4GxxQ5LmKtNs4ABhhCi8paZ2jN14o5XEQw5JCNvQY6pQH0VoHdFUOlAjfFM0uga7ZTbjRCtQDAqvyYoPoADt2DK5u807oOAt19sJak2URaXAWdwy0cuNgxt6NUw0Cki1b%2BEI92foKPDdORActx%2Fd3znU374jIlq2tf9NYhoANBufDORJtbLJIw8rDISVRWemRmkaOzVvFh2v5zrUx6jTkmVoIOiSKh9XVimysJ2UE6z5NefAbjlJcMb4sk1cimS2mO88RSAdg%3D&X-Amz-SignedHeaders=host&response-expires=Wed%2C%2005%20Feb%202025%2023%3A05%3A43%20GMT&X-Amz-Signature=b4709905f423e8389a82e37f26ad448f8b3133a010b41fb926f4af117f7acb08"
$ws.Range("L2").Value = "No video URL"

# Row 3
$ws.Range("A3").Value = "281474991205262-1738703041532"
$ws.Range("B3").Value = "Mobile Usage"
$ws.Range("C3").Value = "2025-02-04T15:04:01.532"
$ws.Range("D3").Value = "281474991205262"
$ws.Range("E3").Value = "132"
$ws.Range("F3").Value = "52215867"
$ws.Range("G3").Value = "EMMANUEL SALCEDO"
$ws.Range("H3").Value = 20.672200369
$ws.Range("I3").Value = -103.29784164
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = "https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1738703039032/BIzNJWuLlc-camera-video-segment-driver-1738703041532.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSK5OXEILC%2F20250205%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250205T150543Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjECsaCXVzLXdlc3QtMiJHMEUCIHsgChGSRKwjEjGoDqjcGh0LrGtLRUqAEs2OolKwO3gIAiEAp4o6%2BVCBhl1VcA%2FXVT6kmUoc4Q%2BBKuNSqGQ4TpPl2m4q3QMIRBAEGgw3ODEyMDQ5NDIyNDQiDI4RT2evNNQhplb4uiq6A9490COmLc8PeJwtKgtCKCn4jpIrTNdV5vXkE6bNQZQffREVTswCw8bTKS3pu5bBuIhAD6V8hA20CSqTdb%2BRUm8yG1Zrie012UZIgvDGLl7YkachS6vz4Yx35THvqPLu%2B9Gmb4gHx3IZW7a9piMFic8LVgMcuPVC647r6O5Ip8SIguX3TDbP0fUliyZB8NqvosYJUBjV5nScUPyAQSz0uR8a9%2FjIOBQbXgiYz07uDubW2VC1aYQoRWQks8T0bMeHny748fSDpEBtnAbGly4BXSMRzR%2BNbIy4CHz1CZYruHs%2BWV5SWee6BPBegyYm1fO2aMY1upjqIn89nbNIFtBkgmpCXRn2H4XAUe1zwU%2FSsptShLhmjIhtaU%2BrIw1CGPirdG5%2FVNxN3NUAMoajto%2Fw6RNKdvA8tspWKn1m5InxeYN164ODNASIzkJE5dHFIX%2B7z%2BJzKC0zbUvXVpQKZbYD6vL6nYoEXBXuI6ZnPetxTIOYfmECP1IlJWR59HLdsCboc0zrS7bz%2BAf5AAaTt14Np4KZKFR9G0Phzo4j9EVle1yi7E4GxxQ5LmKtNs4ABhhCi8paZ2jN14o5XEQw5JCNvQY6pQH0VoHdFUOlAjfFM0uga7ZTbjRCtQDAqvyYoPoADt2DK5u807oOAt19sJak2URaXAWdwy0cuNgxt6NUw0Cki1b%2BEI92foKPDdORActx%2Fd3znU374jIlq2tf9NYhoANBufDORJtbLJIw8rDISVRWemRmkaOzVvFh2v5zrUx6jTkmVoIOiSKh9XVimysJ2UE6z5NefAbjlJcMb4sk1cimS2mO88RSAdg%3D&X-Amz-SignedHeaders=host&response-expires=Wed%2C%2005%20Feb%202025%2023%3A05%3A43%20GMT&X-Amz-Signature=1a4f1dcdc49adafd7785fb67a01b38e504877652218ef1fa892aafae9bcc6b0a"
$ws.Range("L3").Value = "No video URL"

# Row 4
$ws.Range("A4").Value = "281474992631761-1738653345506"
$ws.Range("B4").Value = "No Seat Belt"
$ws.Range("C4").Value = "2025-02-04T01:15:45.506"
$ws.Range("D4").Value = "281474992631761"
$ws.Range("E4").Value = "118"
$ws.Range("F4").Value = "51834005"
$ws.Range("G4").Value = "LUIS FIDENCIO GALINDO BEAS"
$ws.Range("H4").Value = 20.644768619
$ws.Range("I4").Value = -103.356231789
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = "https://s3.samsara.com/samsara-cvdata/4006124/281474992631761/1738653343006/eUbdcuHBpU-camera-video-segment-driver-1738653345506.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSK5OXEILC%2F20250205%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250205T150543Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjECsaCXVzLXdlc3QtMiJHMEUCIHsgChGSRKwjEjGoDqjcGh0LrGtLRUqAEs2OolKwO3gIAiEAp4o6%2BVCBhl1VcA%2FXVT6kmUoc4Q%2BBKuNSqGQ4TpPl2m4q3QMIRBAEGgw3ODEyMDQ5NDIyNDQiDI4RT2evNNQhplb4uiq6A9490COmLc8PeJwtKgtCKCn4jpIrTNdV5vXkE6bNQZQffREVTswCw8bTKS3pu5bBuIhAD6V8hA20CSqTdb%2BRUm8yG1Zrie012UZIgvDGLl7YkachS6vz4Yx35THvqPLu%2B9Gmb4gHx3IZW7a9piMFic8LVgMcuPVC647r6O5Ip8SIguX3TDbP0fUliyZB8NqvosYJUBjV5nScUPyAQSz0uR8a9%2FjIOBQbXgiYz07uDubW2VC1aYQoRWQks8T0bMeHny748fSDpEBtnAbGly4BXSMRzR%2BNbIy4CHz1CZYruHs%2BWV5SWee6BPBegyYm1fO2aMY1upjqIn89nbNIFtBkgmpCXRn2H4XAUe1zwU%2FSsptShLhmjIhtaU%2BrIw1CGPirdG5%2FVNxN3NUAMoajto%2Fw6RNKdvA8tspWKn1m5InxeYN164ODNASIzkJE5dHFIX%2B7z%2BJzKC0zbUvXVpQKZbYD6vL6nYoEXBXuI6ZnPetxTIOYfmECP1IlJWR59HLdsCboc0zrS7bz%2BAf5AAaTt14Np4KZKFR9G0Phzo4j9EVle1yi7E4GxxQ5LmKtNs4ABhhCi8paZ2jN14o5XEQw5JCNvQY6pQH0VoHdFUOlAjfFM0uga7ZTbjRCtQDAqvyYoPoADt2DK5u807oOAt19sJak2URaXAWdwy0cuNgxt6NUw0Cki1b%2BEI92foKPDdORActx%2Fd3znU374jIlq2tf9NYhoANBufDORJtbLJIw8rDISVRWemRmkaOzVvFh2v5zrUx6jTkmVoIOiSKh9XVimysJ2UE6z5NefAbjlJcMb4sk1cimS2mO88RSAdg%3D&X-Amz-SignedHeaders=host&response-expires=Wed%2C%2005%20Feb%202025%2023%3A05%3A43%20GMT&X-Amz-Signature=7e9ac590ab1cd73cf97cff857952f62e0c6dc58b5de8c10fe83748958d09b887"
$ws.Range("L4").Value = "No video URL"
